$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force text to avoid numeric auto-conversion ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D2").Value = "66.763.65"
$ws.Range("D3").Value = "3.088.73"
$ws.Range("D5").Value = "580.05"
$ws.Range("D6").Value = "167.47"
$ws.Range("D8").Value = "3.086.11"
$ws.Range("D12").Value = "0.483"
$ws.Range("D14").Value = "36.67"
$ws.Range("D16").Value = "3.603.32"
$ws.Range("D17").Value = "66.892.18"
$ws.Range("D18").Value = "7.23"
$ws.Range("D19").Value = "3.091.62"
$ws.Range("D20").Value = "16.26"
$ws.Range("D21").Value = "468.75"
$ws.Range("D24").Value = "83.20"
$ws.Range("D26").Value = "12.83"
$ws.Range("D27").Value = "10.11"
$ws.Range("D29").Value = "8.04"
$ws.Range("D33").Value = "28.21"
$ws.Range("D38").Value = "2.13"
$ws.Range("D39").Value = "46.57"
$ws.Range("D40").Value = "50.26"
$ws.Range("D44").Value = "2.84"
$ws.Range("D46").Value = "384.05"
$ws.Range("D47").Value = "2.769.45"
$ws.Range("D48").Value = "134.85"
$ws.Range("D50").Value = "24.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +3.04%  "
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +2.21%  "
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +4.96%  "
$ws.Range("E9").Value = "  +1.24%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  +2.31%  "
$ws.Range("E12").Value = "  +6.14%  "
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("E14").Value = "  +7.88%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("E16").Value = "  +5.16%  "
$ws.Range("E17").Value = "  +3.03%  "
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("E19").Value = "  +5.18%  "
$ws.Range("E20").Value = "  +17.79%  "
$ws.Range("E21").Value = "  +5.11%  "
$ws.Range("E22").Value = "  +4.82%  "
$ws.Range("E23").Value = "  +4.45%  "
$ws.Range("E24").Value = "  +0.87%  "
$ws.Range("E25").Value = "  +7.02%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +3.83%  "
$ws.Range("E34").Value = "  +4.89%  "
$ws.Range("E35").Value = "  +0.10%  "
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("E38").Value = "  +8.17%  "
$ws.Range("E39").Value = "  +6.49%  "
$ws.Range("E40").Value = "  +2.44%  "
$ws.Range("E41").Value = "  +5.97%  "
$ws.Range("E42").Value = "  +3.43%  "
$ws.Range("E43").Value = "  +3.34%  "
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +6.49%  "
$ws.Range("E51").Value = "  +4.73%  "
